$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("disk images")

# Move the "- ""New SmpBank""" text from G5 to G6, keeping G6 formatted with
# the quote-prefix style (as G5 had), and clear G5's contents while keeping
# its existing style.
$ws.Range("G6").Formula = "'- ""New SmpBank"""
$ws.Range("G5").ClearContents()

# Update the active selection on the sheet view from G16 to F7.
$ws.Range("F7").Select()
